$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark near the top of the document
#    (SECTION IV heading paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Rebuild "Additional chapter notes" as six separate runs with a
#    "_GoBack" bookmark inserted between "N" and "ote", matching the target
#    structure. The whole word (including the trailing "s") is reinserted via
#    InsertXML so the run boundaries are preserved exactly as authored -
#    simple text/range edits cause adjacent same-formatted runs to coalesce.
$target = $d.Content
$target.Find.Execute("Additional chapter notes", $true) | Out-Null
$oldRange = $d.Range($target.Start, $target.End)

$fragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Additional </w:t></w:r><w:r><w:t>C</w:t></w:r><w:r><w:t xml:space="preserve">hapter </w:t></w:r><w:r><w:t>N</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>ote</w:t></w:r><w:r w:rsidRPr="00C6181C"><w:t>s</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint = $d.Range($oldRange.Start, $oldRange.Start)
$insertionPoint.InsertXML($fragment)

# The original run text is now pushed after the newly inserted runs; locate it
# again (absolute offsets do not track the insertion) and delete it.
# MatchCase is required so this does not re-match the just-inserted
# "Additional Chapter Notes" text instead of the original lower-case run.
$staleRange = $d.Content
$staleRange.Find.Execute("Additional chapter notes", $true) | Out-Null
$staleRange.Text = ""
